$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 304.5625
$ws.Range("I28").Value = 222.45454
$ws.Range("K28").Value = 222.45454
$ws.Range("M28").Value = 262.54546

$ws.Range("H131").Value = 4438.077
$ws.Range("I131").Value = 1699.1666
$ws.Range("J131").Value = 6785.7144
$ws.Range("K131").Value = 5097.4998
$ws.Range("L131").Value = 20357.1432
$ws.Range("M131").Value = -57.4997999999996
$ws.Range("N131").Value = -30437.1432

$ws.Range("H132").Value = 4066550
$ws.Range("I132").Value = 996.831
$ws.Range("J132").Value = 30307848
$ws.Range("K132").Value = 2990.493
$ws.Range("L132").Value = 90923544
$ws.Range("M132").Value = -460.4929999999999
$ws.Range("N132").Value = -90928604

$ws.Range("H137").Value = 2011.3438
$ws.Range("I137").Value = 1548.6786
$ws.Range("J137").Value = 5250
$ws.Range("K137").Value = 4646.0358
$ws.Range("L137").Value = 15750
$ws.Range("M137").Value = -2096.0358
$ws.Range("N137").Value = -20850

$ws.Range("H141").Value = 2300.4348
$ws.Range("I141").Value = 1647.579
$ws.Range("J141").Value = 5401.5
$ws.Range("K141").Value = 4942.737
$ws.Range("L141").Value = 16204.5
$ws.Range("M141").Value = 237.2629999999999
$ws.Range("N141").Value = -26564.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 134615.94
$ws.Range("I2").Value = 201243.3
$ws.Range("J2").Value = 1361.2
$ws.Range("K2").Value = 201243.3
$ws.Range("L2").Value = 1361.2
$ws.Range("M2").Value = -201130.3
$ws.Range("N2").Value = -1587.2

$ws.Range("H32").Value = 10029.145
$ws.Range("I32").Value = 7058.169
$ws.Range("J32").Value = 21467.4
$ws.Range("K32").Value = 7058.169
$ws.Range("L32").Value = 21467.4
$ws.Range("M32").Value = -6771.169
$ws.Range("N32").Value = -22041.4

$ws.Range("H74").Value = 9092295
$ws.Range("I74").Value = 1125.5294
$ws.Range("J74").Value = 23811330
$ws.Range("K74").Value = 1125.5294
$ws.Range("L74").Value = 23811330
$ws.Range("M74").Value = -251.5293999999999
$ws.Range("N74").Value = -23813078

$ws.Range("H76").Value = 42288
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 42288
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 42288
$ws.Range("M76").Value = ""
$ws.Range("N76").Value = -42964

$ws.Range("H77").Value = 9092295
$ws.Range("I77").Value = 1125.5294
$ws.Range("J77").Value = 23811330
$ws.Range("K77").Value = 5627.646999999999
$ws.Range("L77").Value = 119056650
$ws.Range("M77").Value = -1259.646999999999
$ws.Range("N77").Value = -119065386

$ws.Range("H79").Value = 42288
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 42288
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 42288
$ws.Range("M79").Value = ""
$ws.Range("N79").Value = -44628

$ws.Range("H116").Value = 134615.94
$ws.Range("I116").Value = 201243.3
$ws.Range("J116").Value = 1361.2
$ws.Range("K116").Value = 201243.3
$ws.Range("L116").Value = 1361.2
$ws.Range("M116").Value = -198949.3
$ws.Range("N116").Value = -5949.2

$ws.Range("H122").Value = 1071812.2
$ws.Range("I122").Value = 1168977
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 3506931
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -3504481
$ws.Range("N122").Value = -13900

$ws.Range("H132").Value = 1517726.4
$ws.Range("I132").Value = 1849.4762
$ws.Range("J132").Value = 4170511
$ws.Range("K132").Value = 5548.4286
$ws.Range("L132").Value = 12511533
$ws.Range("M132").Value = -3018.4286
$ws.Range("N132").Value = -12516593

$ws.Range("H138").Value = 28554
$ws.Range("J138").Value = 28554
$ws.Range("L138").Value = 28554
$ws.Range("N138").Value = -38834

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 134615.94
$ws.Range("I3").Value = 201243.3
$ws.Range("J3").Value = 1361.2
$ws.Range("K3").Value = 201243.3
$ws.Range("L3").Value = 1361.2
$ws.Range("M3").Value = -201129.3
$ws.Range("N3").Value = -1589.2

$ws.Range("H94").Value = 1722.1765
$ws.Range("I94").Value = 1397.4445
$ws.Range("J94").Value = 2087.5
$ws.Range("K94").Value = 1397.4445
$ws.Range("L94").Value = 2087.5
$ws.Range("M94").Value = -946.4445000000001
$ws.Range("N94").Value = -2989.5

$ws.Range("H107").Value = 202231.33
$ws.Range("I107").Value = 216547.86
$ws.Range("K107").Value = 216547.86
$ws.Range("M107").Value = -214627.86

$ws.Range("H134").Value = 33389.027
$ws.Range("I134").Value = 6351.6206
$ws.Range("J134").Value = 145401.14
$ws.Range("K134").Value = 19054.8618
$ws.Range("L134").Value = 436203.42
$ws.Range("M134").Value = -16519.8618
$ws.Range("N134").Value = -441273.42

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1580.2646
$ws.Range("I16").Value = 1456.48
$ws.Range("J16").Value = 1924.1111
$ws.Range("K16").Value = 1456.48
$ws.Range("L16").Value = 1924.1111
$ws.Range("M16").Value = -1169.48
$ws.Range("N16").Value = -2498.1111

$ws.Range("H31").Value = 11719231
$ws.Range("I31").Value = 2031.35
$ws.Range("J31").Value = 21483564
$ws.Range("K31").Value = 2031.35
$ws.Range("L31").Value = 21483564
$ws.Range("M31").Value = -1736.35
$ws.Range("N31").Value = -21484154

$ws.Range("H34").Value = 11719231
$ws.Range("I34").Value = 2031.35
$ws.Range("J34").Value = 21483564
$ws.Range("K34").Value = 2031.35
$ws.Range("L34").Value = 21483564
$ws.Range("M34").Value = -1829.35
$ws.Range("N34").Value = -21483968

$ws.Range("H58").Value = 5522731.5
$ws.Range("I58").Value = 9805698
$ws.Range("J58").Value = 668702.9399999999
$ws.Range("K58").Value = 9805698
$ws.Range("L58").Value = 668702.9399999999
$ws.Range("M58").Value = -9805495
$ws.Range("N58").Value = -669108.9399999999

$ws.Range("H113").Value = 1580.2646
$ws.Range("I113").Value = 1456.48
$ws.Range("J113").Value = 1924.1111
$ws.Range("K113").Value = 1456.48
$ws.Range("L113").Value = 1924.1111
$ws.Range("M113").Value = 713.52
$ws.Range("N113").Value = -6264.1111

$ws.Range("H132").Value = 4168700
$ws.Range("I132").Value = 5883579
$ws.Range("J132").Value = 3994.2856
$ws.Range("K132").Value = 17650737
$ws.Range("L132").Value = 11982.8568
$ws.Range("M132").Value = -17648207
$ws.Range("N132").Value = -17042.8568

$ws.Range("H136").Value = 5522731.5
$ws.Range("I136").Value = 9805698
$ws.Range("J136").Value = 668702.9399999999
$ws.Range("K136").Value = 29417094
$ws.Range("L136").Value = 2006108.82
$ws.Range("M136").Value = -29414544
$ws.Range("N136").Value = -2011208.82

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 3252011
$ws.Range("I5").Value = 7690.2856
$ws.Range("J5").Value = 6496331.5
$ws.Range("K5").Value = 23070.8568
$ws.Range("L5").Value = 19488994.5
$ws.Range("M5").Value = -22958.8568
$ws.Range("N5").Value = -19489218.5

$ws.Range("H34").Value = 2611.111
$ws.Range("I34").Value = 1200
$ws.Range("J34").Value = 2893.3333
$ws.Range("K34").Value = 3600
$ws.Range("L34").Value = 8679.999899999999
$ws.Range("M34").Value = -3516
$ws.Range("N34").Value = -8847.999899999999

$ws.Range("H39").Value = 2133.1667
$ws.Range("J39").Value = 2133.1667
$ws.Range("L39").Value = 6399.500100000001
$ws.Range("N39").Value = -6987.500100000001

$ws.Range("H55").Value = 6725
$ws.Range("J55").Value = 6725
$ws.Range("L55").Value = 20175
$ws.Range("N55").Value = -20529

$ws.Range("H101").Value = 8266.666999999999
$ws.Range("J101").Value = 8266.666999999999
$ws.Range("L101").Value = 24800.001
$ws.Range("N101").Value = -29668.001

$ws.Range("H107").Value = 442.39285
$ws.Range("I107").Value = 400.29413
$ws.Range("J107").Value = 507.45456
$ws.Range("K107").Value = 1200.88239
$ws.Range("L107").Value = 1522.36368
$ws.Range("M107").Value = 719.11761
$ws.Range("N107").Value = -5362.36368

$ws.Range("H135").Value = 3252011
$ws.Range("I135").Value = 7690.2856
$ws.Range("J135").Value = 6496331.5
$ws.Range("K135").Value = 69212.5704
$ws.Range("L135").Value = 58466983.5
$ws.Range("M135").Value = -66677.5704
$ws.Range("N135").Value = -58472053.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2975.125
$ws.Range("I102").Value = 2845.6365
$ws.Range("J102").Value = 3260
$ws.Range("K102").Value = 2845.6365
$ws.Range("L102").Value = 3260
$ws.Range("M102").Value = -1223.6365
$ws.Range("N102").Value = -6504

$ws.Range("H122").Value = 48745540
$ws.Range("I122").Value = 56044556
$ws.Range("J122").Value = 33336500
$ws.Range("K122").Value = 168133668
$ws.Range("L122").Value = 100009500
$ws.Range("M122").Value = -168131218
$ws.Range("N122").Value = -100014400

$ws.Range("H132").Value = 5213274
$ws.Range("I132").Value = 10422405
$ws.Range("K132").Value = 31267215
$ws.Range("M132").Value = -31264685

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 438.5263
$ws.Range("I16").Value = 451.30768
$ws.Range("J16").Value = 410.83334
$ws.Range("K16").Value = 451.30768
$ws.Range("L16").Value = 410.83334
$ws.Range("M16").Value = -281.30768
$ws.Range("N16").Value = -750.83334

$ws.Range("H41").Value = 11500
$ws.Range("J41").Value = 11500
$ws.Range("L41").Value = 11500
$ws.Range("N41").Value = -12376

$ws.Range("H68").Value = 47621110
$ws.Range("I68").Value = 1938.2
$ws.Range("J68").Value = 90911270
$ws.Range("K68").Value = 1938.2
$ws.Range("L68").Value = 90911270
$ws.Range("M68").Value = -1189.2
$ws.Range("N68").Value = -90912768

$ws.Range("H71").Value = 47621110
$ws.Range("I71").Value = 1938.2
$ws.Range("J71").Value = 90911270
$ws.Range("K71").Value = 9691
$ws.Range("L71").Value = 454556350
$ws.Range("M71").Value = -5947
$ws.Range("N71").Value = -454563838

$ws.Range("H100").Value = 1896
$ws.Range("I100").Value = 1802
$ws.Range("J100").Value = 1990
$ws.Range("K100").Value = 1802
$ws.Range("L100").Value = 1990
$ws.Range("M100").Value = -1261
$ws.Range("N100").Value = -3072

$ws.Range("H132").Value = 21747678
$ws.Range("I132").Value = 23818366
$ws.Range("K132").Value = 71455098
$ws.Range("M132").Value = -71452568

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H104").Value = 40728.332
$ws.Range("J104").Value = 40728.332
$ws.Range("L104").Value = 40728.332
$ws.Range("N104").Value = -47716.332

$ws.Range("H126").Value = 1465.2222
$ws.Range("I126").Value = 1439.4
$ws.Range("J126").Value = 1497.5
$ws.Range("K126").Value = 4318.200000000001
$ws.Range("L126").Value = 4492.5
$ws.Range("M126").Value = -1848.200000000001
$ws.Range("N126").Value = -9432.5

$ws.Range("H132").Value = 2883.6667
$ws.Range("I132").Value = 1238.5
$ws.Range("K132").Value = 3715.5
$ws.Range("M132").Value = -1185.5

$ws.Range("H138").Value = 52114.5
$ws.Range("J138").Value = 54229
$ws.Range("L138").Value = 54229
$ws.Range("N138").Value = -64509
